$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 1593
$ws1.Range("F10").Value = 2656
$ws1.Range("F11").Value = 2656
$ws1.Range("F13").Value = 1721
$ws1.Range("F15").Value = 260
$ws1.Range("F16").Value = 674
$ws1.Range("F17").Value = 4916
$ws1.Range("F18").Value = 151
$ws1.Range("F21").Value = 3381
$ws1.Range("F22").Value = 850
$ws1.Range("F25").Value = 34
$ws1.Range("F26").Value = 2402
$ws1.Range("F32").Value = 1278
$ws1.Range("F35").Value = 17
$ws1.Range("F37").Value = 1385
$ws1.Range("F38").Value = 1359

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 127
$ws2.Range("F18").Value = 255
$ws2.Range("F19").Value = 514

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 832
$ws3.Range("F4").Value = 236

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 832
$ws4.Range("F7").Value = 236
$ws4.Range("F14").Value = 1593
$ws4.Range("F18").Value = 2656
$ws4.Range("F21").Value = 1721
$ws4.Range("F22").Value = 127
$ws4.Range("F24").Value = 260
$ws4.Range("F25").Value = 674
$ws4.Range("F26").Value = 4916
$ws4.Range("F29").Value = 3381
$ws4.Range("F33").Value = 34
$ws4.Range("F34").Value = 2402
$ws4.Range("F38").Value = 1278
$ws4.Range("F40").Value = 255
$ws4.Range("F41").Value = 514
$ws4.Range("F45").Value = 17
$ws4.Range("F47").Value = 1385
$ws4.Range("F49").Value = 1359
